# Natmi following Dr Hou advice
# Replace the single FAPs->FAPs communication row with the full set of
# pairwise ECs/FAPs/sCs -> ECs/FAPs/sCs (minus the ECs target rows)
# Agt-Agtr2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Agt", "Agtr2", "FAPs", 1, 0.3333333333333333, 0.02619133333333333, 0.078574,   0.00783286589743605, 0.007832865897436052, 3, 1,                  2.002414333333333, 6.007243, 0.8176262304899502, 0.8176262304899502, 0.05244590127577778, 0.472013111482,   0.006404356617653918, 0.00640435661765392),
    @("ECs",  "Agt", "Agtr2", "sCs",  1, 0.3333333333333333, 0.02619133333333333, 0.078574,   0.00783286589743605, 0.007832865897436052, 2, 0.6666666666666666, 0.446644,           1.339932, 0.1823737695100498, 0.1823737695100498, 0.01169820188533333, 0.105283816968,   0.001428509279782132, 0.001428509279782132),
    @("FAPs", "Agt", "Agtr2", "FAPs", 3, 1,                  3.276419666666667,   9.829259,   0.9798567925543612,  0.9798567925543613,  3, 1,                  2.002414333333333, 6.007243, 0.8176262304899502, 0.8176262304899502, 6.560749702548556,   59.046747322937,  0.8011566157161955,  0.8011566157161956),
    @("FAPs", "Agt", "Agtr2", "sCs",  3, 1,                  3.276419666666667,   9.829259,   0.9798567925543612,  0.9798567925543613,  2, 0.6666666666666666, 0.446644,           1.339932, 0.1823737695100498, 0.1823737695100498, 1.463393185598667,   13.170538670388,  0.1787001768381657,  0.1787001768381658),
    @("sCs",  "Agt", "Agtr2", "FAPs", 1, 0.3333333333333333, 0.041163,             0.123489,   0.01231034154820272, 0.01231034154820272, 3, 1,                  2.002414333333333, 6.007243, 0.8176262304899502, 0.8176262304899502, 0.082425381203,      0.741828430827,   0.01006525815610081, 0.01006525815610081),
    @("sCs",  "Agt", "Agtr2", "sCs",  1, 0.3333333333333333, 0.041163,             0.123489,   0.01231034154820272, 0.01231034154820272, 2, 0.6666666666666666, 0.446644,           1.339932, 0.1823737695100498, 0.1823737695100498, 0.018385206972,      0.165466862748,   0.002245083392101912, 0.002245083392101912)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
